$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (e.g. "1.000" -> 1). Pre-format as Text so the literal string is preserved,
# matching the source workbook where these are inline/shared strings.
$textCells = @('D4', 'D5', 'D6', 'D8', 'D9', 'D10', 'D11', 'D13', 'D14', 'D15', 'D17', 'D18', 'D19', 'D20', 'D21', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D33', 'D34', 'D36', 'D37', 'D39', 'D41', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price / Volume(1h) values scraped for this run.
$ws.Range('D2').Value = '29.217.44'
$ws.Range('E2').Value = '  +0.49%  '
$ws.Range('D3').Value = '1.859.60'
$ws.Range('E3').Value = '  +1.01%  '
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '0.7021'
$ws.Range('E5').Value = '  +1.46%  '
$ws.Range('D6').Value = '237.34'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '0.07737'
$ws.Range('E8').Value = '  +4.51%  '
$ws.Range('D9').Value = '0.3045'
$ws.Range('E9').Value = '  +0.81%  '
$ws.Range('D10').Value = '23.26'
$ws.Range('E10').Value = '  +0.48%  '
$ws.Range('D11').Value = '0.08174'
$ws.Range('E11').Value = '  +1.20%  '
$ws.Range('D12').Value = '1.864.25'
$ws.Range('E12').Value = '  -0.71%  '
$ws.Range('D13').Value = '0.7185'
$ws.Range('E13').Value = '  -0.01%  '
$ws.Range('D14').Value = '5.159'
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('D15').Value = '89.11'
$ws.Range('E15').Value = '  +0.56%  '
$ws.Range('D16').Value = '29.225.47'
$ws.Range('E16').Value = '  +0.88%  '
$ws.Range('D17').Value = '5.767'
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('D18').Value = '13.35'
$ws.Range('E18').Value = '  +3.34%  '
$ws.Range('D19').Value = '0.000007733'
$ws.Range('E19').Value = '  +1.47%  '
$ws.Range('D20').Value = '236.84'
$ws.Range('E20').Value = '  -1.22%  '
$ws.Range('D21').Value = '0.9998'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').Value = '2.108.57'
$ws.Range('E22').Value = '  +2.24%  '
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').Value = '7.445'
$ws.Range('E24').Value = '  -1.85%  '
$ws.Range('D25').Value = '0.1483'
$ws.Range('E25').Value = '  +1.40%  '
$ws.Range('D26').Value = '162.27'
$ws.Range('E26').Value = '  +0.57%  '
$ws.Range('D27').Value = '9.005'
$ws.Range('E27').Value = '  +0.49%  '
$ws.Range('D28').Value = '18.00'
$ws.Range('E28').Value = '  +0.24%  '
$ws.Range('D29').Value = '2.039'
$ws.Range('E29').Value = '  +6.44%  '
$ws.Range('D30').Value = '1.432'
$ws.Range('E30').Value = '  +4.19%  '
$ws.Range('D31').Value = '4.434'
$ws.Range('E31').Value = '  +0.25%  '
$ws.Range('E32').Value = '  -0.23%  '
$ws.Range('D33').Value = '4.031'
$ws.Range('E33').Value = '  +0.65%  '
$ws.Range('D34').Value = '0.05225'
$ws.Range('E34').Value = '  +0.90%  '
$ws.Range('E35').Value = '  -0.70%  '
$ws.Range('D36').Value = '0.7074'
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').Value = '1.000'
$ws.Range('E37').Value = '  +0.34%  '
$ws.Range('E38').Value = '  +0.47%  '
$ws.Range('D39').Value = '0.01845'
$ws.Range('E39').Value = '  -0.76%  '
$ws.Range('E40').Value = '  +1.78%  '
$ws.Range('D41').Value = '0.9350'
$ws.Range('E41').Value = '  +2.83%  '
$ws.Range('D42').Value = '1.141.77'
$ws.Range('E42').Value = '  +8.34%  '
$ws.Range('D43').Value = '0.4273'
$ws.Range('E43').Value = '  +0.24%  '
$ws.Range('D44').Value = '5.901'
$ws.Range('E44').Value = '  +0.30%  '
$ws.Range('D45').Value = '70.79'
$ws.Range('E45').Value = '  +1.90%  '
$ws.Range('D46').Value = '1.0000'
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').Value = '103.34'
$ws.Range('E47').Value = '  +1.00%  '
$ws.Range('D48').Value = '1.794'
$ws.Range('E48').Value = '  +3.44%  '
$ws.Range('D49').Value = '2.004.83'
$ws.Range('E49').Value = '  +2.33%  '
$ws.Range('D50').Value = '9.173'
$ws.Range('E50').Value = '  -0.11%  '
$ws.Range('D51').Value = '6.963'
$ws.Range('E51').Value = '  -1.98%  '

# Restore default (General) styling now that the literal text is committed,
# so these cells don't end up carrying an explicit style index.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
